# MHD2-259: Report template and related changes for reporting on 136 genes
#
# 1) Fix typo "DELECTED" -> "DETECTED" in the "NO VARIANTS DETECTED." summary text.
# 2) Trim the CDS footnote paragraph down to just "CDS – coding sequence"
#    (drop the trailing "; * - partial coverage; ex – exon; int - intron;
#    please note FLT3-ITDs and UBTF-TDs are not reliably detected with this
#    assay. A separate assay may have been performed, result included in
#    Test Description if sample tested." text).
# 3) Update the cached SAVEDATE field result from "31-Oct-2025" to "7-Nov-2025".

$d = $word.ActiveDocument

# --- 1) Typo fix: DELECTED -> DETECTED -------------------------------------
$d.Content.Find.Execute("DELECTED", $true, $false, $false, $false, $false, `
    $true, 1, $false, "DETECTED", 2) | Out-Null

# --- 2) Trim the CDS footnote paragraph -------------------------------------
$cdsStart = $d.Content
$cdsStart.Find.Execute("; * - partial coverage; ex", $true, $false, $false, `
    $false, $false, $true, 1, $false, "", 0) | Out-Null
$startPos = $cdsStart.Start

$cdsEnd = $d.Content
$cdsEnd.Find.Execute("if sample tested.", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null
$endPos = $cdsEnd.End

$trimRange = $d.Range($startPos, $endPos)
$trimRange.Text = ""

# --- 3) Update the cached SAVEDATE field text -------------------------------
$d.Content.Find.Execute("31-Oct-2025", $true, $false, $false, $false, $false, `
    $true, 1, $false, "7-Nov-2025", 2) | Out-Null
